# Foglio1 ("test nazione" sheet) - aggiungo colonna F "nazione" con i valori
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# I valori vengono inseriti prima nella colonna dati, poi l'intestazione,
# cosi' come e' stato fatto nel file originale (ordine di inserimento
# delle shared strings: Italia, Calimera, nazione).
$ws.Range("F2").Value = "Italia"
$ws.Range("F3").Value = "Calimera"
$ws.Range("F1").Value = "nazione"

$ws.Range("F1").Select()
